$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had 12 data rows (4 sending clusters x 3 target clusters,
# including "ECs" as a target). The refreshed TPM run drops "ECs" as a target
# cluster, so only 2 target clusters remain -> 4 x 2 = 8 data rows (rows 2-9).
# Delete the now-unused trailing rows (previously rows 10-13) first so the used
# range / dimension shrinks to A1:T9, then rewrite rows 2-9 with the recomputed values.
$ws.Rows("10:13").Delete()

$newData = @(
    @("ECs", "Tnfrsf14", "Btla", "MuSCs", 3, 1, 5.209944, 15.629832, 0.1751928672265232, 0.1751928672265232, 1, 0.3333333333333333, 0.09278333333333333, 0.27835, 0.07826677681470105, 0.07826677681470104, 0.4833959708, 4.3505637372, 0.01371178103874585, 0.01371178103874584),
    @("ECs", "Tnfrsf14", "Btla", "Resolving-Mac", 3, 1, 5.209944, 15.629832, 0.1751928672265232, 0.1751928672265232, 3, 1, 1.092692, 3.278076, 0.921733223185299, 0.9217332231852989, 5.692864129248, 51.235777163232, 0.1614810861877774, 0.1614810861877773),
    @("FAPs", "Tnfrsf14", "Btla", "MuSCs", 3, 1, 6.497702, 19.493106, 0.2184958310038485, 0.2184958310038485, 1, 0.3333333333333333, 0.09278333333333333, 0.27835, 0.07826677681470105, 0.07826677681470104, 0.6028784505666667, 5.4259060551, 0.01710096444012085, 0.01710096444012084),
    @("FAPs", "Tnfrsf14", "Btla", "Resolving-Mac", 3, 1, 6.497702, 19.493106, 0.2184958310038485, 0.2184958310038485, 3, 1, 1.092692, 3.278076, 0.921733223185299, 0.9217332231852989, 7.099986993784, 63.899882944056, 0.2013948665637276, 0.2013948665637276),
    @("MuSCs", "Tnfrsf14", "Btla", "MuSCs", 3, 1, 1.714656666666667, 5.14397, 0.057658127945791, 0.057658127945791, 1, 0.3333333333333333, 0.09278333333333333, 0.27835, 0.07826677681470105, 0.07826677681470104, 0.1590915610555556, 1.4318240495, 0.004512715831486702, 0.004512715831486701),
    @("MuSCs", "Tnfrsf14", "Btla", "Resolving-Mac", 3, 1, 1.714656666666667, 5.14397, 0.057658127945791, 0.057658127945791, 3, 1, 1.092692, 3.278076, 0.921733223185299, 0.9217332231852989, 1.873591622413333, 16.86232460172, 0.05314541211430431, 0.0531454121143043),
    @("Resolving-Mac", "Tnfrsf14", "Btla", "MuSCs", 3, 1, 16.31603133333333, 48.948094, 0.5486531738238374, 0.5486531738238374, 1, 0.3333333333333333, 0.09278333333333333, 0.27835, 0.07826677681470105, 0.07826677681470104, 1.513855773877778, 13.6247019649, 0.04294131550434766, 0.04294131550434765),
    @("Resolving-Mac", "Tnfrsf14", "Btla", "Resolving-Mac", 3, 1, 16.31603133333333, 48.948094, 0.5486531738238374, 0.5486531738238374, 3, 1, 1.092692, 3.278076, 0.921733223185299, 0.9217332231852989, 17.82839690968266, 160.455572187144, 0.5057118583194897, 0.5057118583194897)
)

$startRow = 2
for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = $newData[$i]
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($startRow + $i, $c).Value = $row[$c - 1]
    }
}